# Fruta / hortaliza, semanal
# Insert a new weekly record at row 468 for "Macroferia Regional de Talca - Zanahoria".
# Inserting the row shifts all existing rows 468..532 down to 469..533 (preserving
# their original data), matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 468, pushing the rest of the table down.
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A468").Value2 = 5
$ws.Range("B468").Value2 = "Macroferia Regional de Talca"
$ws.Range("C468").Value2 = "Maule"
$ws.Range("D468").Value2 = 45077
$ws.Range("E468").Value2 = 7
$ws.Range("F468").Value2 = 100114013
$ws.Range("G468").Value2 = "Zanahoria"
$ws.Range("H468").Value2 = "Sin especificar"
$ws.Range("I468").Value2 = "Primera"
$ws.Range("J468").Value2 = 600
$ws.Range("K468").Value2 = 5000
$ws.Range("L468").Value2 = 5000
$ws.Range("M468").Value2 = 5000
$ws.Range("N468").Value2 = "`$/saco 20 kilos"
$ws.Range("O468").Value2 = "Región de Ñuble"
$ws.Range("P468").Value2 = 250
$ws.Range("Q468").Value2 = 20
$ws.Range("R468").Value2 = "Hortaliza"
